$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text entry so the dd-mm-yyyy strings are kept as literal text
# (matching shared strings) instead of being auto-converted to dates.
$ws.Range("A6:A7").NumberFormat = "@"

# Row 6: 04-08-2021
$ws.Range("A6").Value = "04-08-2021"
$ws.Range("B6").Value = 10000
$ws.Range("C6").Value = 17000
$ws.Range("D6").Value = 5000
$ws.Range("E6").Value = 5000
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 1.39

# Row 7: 05-08-2021
$ws.Range("A7").Value = "05-08-2021"
$ws.Range("B7").Value = 10000
$ws.Range("C7").Value = 16000
$ws.Range("D7").Value = 5000
$ws.Range("E7").Value = 3000
$ws.Range("F7").Value = 2000
$ws.Range("G7").Value = 1.43

# Restore default (unstyled) appearance to match the rest of column A.
$ws.Range("A6:A7").Style = "Normal"
